# Auto-generated Excel COM-interop edit script
# Applies the cfb_weather.xlsx update: refresh weather/odds data and bump Timestamp to 2024-12-02T16:21:13.498323

$wb = $excel.ActiveWorkbook

# ---- Sheet "FBS" (first sheet): rows 2 and 3 swap games and refresh values ----
$ws1 = $wb.Worksheets.Item("FBS")

# Row 2 -> UNLV @ Boise State
$ws1.Range("A2").Value = "UNLV @ Boise State"
$ws1.Range("B2").Value = "FRI 12/06"
$ws1.Range("C2").Value = "06:00 PM"
$ws1.Range("D2").Value = "High"
$ws1.Range("E2").Value = "N-S"
$ws1.Range("F2").Value = "Med"
$ws1.Range("G2").Value = "E"
$ws1.Range("I2").Value = 53.65
$ws1.Range("J2").Value = 70.04
$ws1.Range("K2").Value = 6.8
$ws1.Range("L2").Value = 1970
$ws1.Range("M2").Value = "SSE"
$ws1.Range("N2").Value = "SE"
$ws1.Range("O2").Value = 42.98
$ws1.Range("P2").Value = 3.8
$ws1.Range("Q2").Value = "SE"
$ws1.Range("R2").Value = 0
$ws1.Range("S2").Value = 0
$ws1.Range("T2").Value = 0
$ws1.Range("U2").Value = -3
$ws1.Range("V2").Value = "43.6028839, -116.1958882"
$ws1.Range("W2").Value = 58.5
$ws1.Range("X2").Value = -110
$ws1.Range("Y2").Value = 58.5
$ws1.Range("Z2").Value = -110
$ws1.Range("AA2").Value = -4
$ws1.Range("AB2").Value = -4
$ws1.Range("AE2").Value = 0
$ws1.Range("AF2").Value = 0
$ws1.Range("AK2").Value = "2024-12-02T16:21:13.498323"
$ws1.Range("H2").ClearContents()

# Row 3 -> Marshall @ Louisiana
$ws1.Range("A3").Value = "Marshall @ Louisiana"
$ws1.Range("B3").Value = "SAT 12/07"
$ws1.Range("C3").Value = "06:30 PM"
$ws1.Range("D3").Value = "Low"
$ws1.Range("E3").Value = "NW-SE"
$ws1.Range("F3").Value = "High"
$ws1.Range("H3").Value = -160.176310297
$ws1.Range("I3").Value = 69.21
$ws1.Range("J3").Value = 57.18
$ws1.Range("K3").Value = 9.9
$ws1.Range("L3").Value = 1971
$ws1.Range("M3").Value = "NNW"
$ws1.Range("N3").Value = "NNW"
$ws1.Range("O3").Value = 64.85
$ws1.Range("P3").Value = 7.9
$ws1.Range("Q3").Value = "NNW"
$ws1.Range("R3").Value = 0
$ws1.Range("S3").Value = 0
$ws1.Range("T3").Value = 0
$ws1.Range("U3").Value = -2
$ws1.Range("V3").Value = "30.2158434, -92.0417371"
$ws1.Range("W3").Value = 58.5
$ws1.Range("X3").Value = -110
$ws1.Range("Y3").Value = 56.5
$ws1.Range("Z3").Value = -105
$ws1.Range("AA3").Value = -3.5
$ws1.Range("AB3").Value = -5
$ws1.Range("AE3").Value = -0.03418803418803419
$ws1.Range("AF3").Value = 1.5
$ws1.Range("AK3").Value = "2024-12-02T16:21:13.498323"
$ws1.Range("G3").ClearContents()

# ---- Sheet "Other" (second sheet): refresh values for all 4 game rows ----
$ws2 = $wb.Worksheets.Item("Other")

# Row 1
$ws2.Range("B1").Value = "Home Team"
$ws2.Range("C1").Value = "Away Team"

# Row 2
$ws2.Range("A2").Value = "Villanova vs Incarnate Word"
$ws2.Range("B2").Value = "Incarnate Word"
$ws2.Range("C2").Value = "Villanova"
$ws2.Range("D2").Value = "SAT 12/07"
$ws2.Range("E2").Value = "01:00 PM"
$ws2.Range("F2").Value = "Low"
$ws2.Range("J2").Value = 81.08228299999999
$ws2.Range("K2").Value = 70.74
$ws2.Range("L2").Value = 55.05
$ws2.Range("N2").Value = 2008
$ws2.Range("O2").Value = "NE"
$ws2.Range("P2").Value = "ESE"
$ws2.Range("Q2").Value = 58.28
$ws2.Range("R2").Value = 4.2
$ws2.Range("S2").Value = "ESE"
$ws2.Range("T2").Value = 2.4
$ws2.Range("U2").Value = -1.5
$ws2.Range("V2").Value = 0
$ws2.Range("X2").Value = "29.4674787, -98.470014"

# Row 3
$ws2.Range("A3").Value = "Rhode Island vs Mercer"
$ws2.Range("B3").Value = "Mercer"
$ws2.Range("C3").Value = "Rhode Island"
$ws2.Range("D3").Value = "SAT 12/07"
$ws2.Range("E3").Value = "02:00 PM"
$ws2.Range("F3").Value = "Low"
$ws2.Range("J3").Value = 105.98195272
$ws2.Range("K3").Value = 64.83
$ws2.Range("L3").Value = 52.81
$ws2.Range("N3").Value = 2013
$ws2.Range("O3").Value = "N"
$ws2.Range("P3").Value = "NNE"
$ws2.Range("Q3").Value = 57.62
$ws2.Range("R3").Value = 3
$ws2.Range("S3").Value = "NNE"
$ws2.Range("T3").Value = 0
$ws2.Range("U3").Value = 0
$ws2.Range("V3").Value = 0
$ws2.Range("X3").Value = "32.8262075, -83.6522485"

# Row 4
$ws2.Range("A4").Value = "Montana vs South Dakota State"
$ws2.Range("B4").Value = "South Dakota State"
$ws2.Range("C4").Value = "Montana"
$ws2.Range("D4").Value = "SAT 12/07"
$ws2.Range("E4").Value = "01:00 PM"
$ws2.Range("F4").Value = "Mid"
$ws2.Range("J4").Value = -474.5684815
$ws2.Range("K4").Value = 46.7
$ws2.Range("L4").Value = 47.64
$ws2.Range("N4").Value = 2016
$ws2.Range("O4").Value = "SSE"
$ws2.Range("P4").Value = "SSE"
$ws2.Range("Q4").Value = 35.3
$ws2.Range("R4").Value = 3.3
$ws2.Range("S4").Value = "SSE"
$ws2.Range("T4").Value = 0
$ws2.Range("U4").Value = 0
$ws2.Range("V4").Value = 0
$ws2.Range("X4").Value = "44.3210182, -96.7801386"

# Row 5
$ws2.Range("A5").Value = "Illinois State vs UC Davis"
$ws2.Range("B5").Value = "UC Davis"
$ws2.Range("C5").Value = "Illinois State"
$ws2.Range("D5").Value = "SAT 12/07"
$ws2.Range("E5").Value = "01:00 PM"
$ws2.Range("F5").Value = "High"
$ws2.Range("J5").Value = -231.4896765
$ws2.Range("K5").Value = 62.21
$ws2.Range("L5").Value = 53.1
$ws2.Range("N5").Value = 2007
$ws2.Range("O5").Value = "SE"
$ws2.Range("P5").Value = "SSE"
$ws2.Range("Q5").Value = 62.48
$ws2.Range("R5").Value = 3.3
$ws2.Range("S5").Value = "SSE"
$ws2.Range("T5").Value = 0
$ws2.Range("U5").Value = 0
$ws2.Range("V5").Value = 0
$ws2.Range("X5").Value = "38.5365266, -121.7627936"

